# Automatische test-sync: 2025-08-04 20:20:50
# Adds a new log row (#6, EcoPro-700 stock question) to the "Logs" sheet,
# extends the conditional-formatting ranges to cover it, and bumps the
# "Inkoop / Bestellingen" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append new row 12 to the Logs sheet -------------------------------
$newRow = 12
$logs.Cells.Item($newRow, 1).Value = "Hebben we EcoPro-700 nog op voorraad?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #6: Hebben we EcoPro-700 nog op voorraad?"
$logs.Cells.Item($newRow, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-04 20:20:43"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- Extend conditional formatting ranges from row 11 to row 12 --------
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $col + "2:" + $col + "11"
    $newRange = $col + "2:" + $col + "12"
    $fcs = $logs.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

# --- Update the Dashboard summary count for "Inkoop / Bestellingen" ----
$dashboard.Range("B5").Value = 2
